# Apply cryptos list price/volume updates (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D/E are plain text in this sheet. Force a Text number format on the
# D cells whose new values would otherwise be auto-parsed as numbers by Excel,
# so they remain stored as text exactly as authored.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.092.92"
$ws.Range("E2").Value = "  -0.39%  "
$ws.Range("D3").Value = "1.628.12"
$ws.Range("E3").Value = "  -1.18%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "215.87"
$ws.Range("E5").Value = "  -1.27%  "
$ws.Range("D6").Value = "0.514"
$ws.Range("E6").Value = "  +0.74%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "0.252"
$ws.Range("E8").Value = "  -1.60%  "
$ws.Range("D9").Value = "0.0622"
$ws.Range("E9").Value = "  -0.92%  "
$ws.Range("D10").Value = "19.99"
$ws.Range("E10").Value = "  -0.44%  "
$ws.Range("D11").Value = "0.0850"
$ws.Range("D12").Value = "1.857.50"
$ws.Range("E12").Value = "  -1.09%  "
$ws.Range("D13").Value = "1.614.56"
$ws.Range("E13").Value = "  -0.50%  "
$ws.Range("D14").Value = "4.10"
$ws.Range("E14").Value = "  -0.79%  "
$ws.Range("D15").Value = "0.537"
$ws.Range("E15").Value = "  -0.27%  "
$ws.Range("D16").Value = "64.88"
$ws.Range("E16").Value = "  -3.86%  "
$ws.Range("D17").Value = "27.073.57"
$ws.Range("E17").Value = "  -0.39%  "
$ws.Range("D18").Value = "0.0₃0729"
$ws.Range("E18").Value = "  -1.58%  "
$ws.Range("D19").Value = "213.10"
$ws.Range("E19").Value = "  -3.09%  "
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("D21").Value = "6.78"
$ws.Range("E21").Value = "  +0.14%  "
$ws.Range("D22").Value = "4.36"
$ws.Range("E22").Value = "  -1.83%  "
$ws.Range("D23").Value = "2.47"
$ws.Range("E23").Value = "  -2.05%  "
$ws.Range("D24").Value = "9.06"
$ws.Range("E24").Value = "  -1.78%  "
$ws.Range("D25").Value = "147.14"
$ws.Range("E25").Value = "  -0.76%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").Value = "7.29"
$ws.Range("E27").Value = "  -1.51%  "
$ws.Range("E28").Value = "  -1.24%  "
$ws.Range("D29").Value = "15.52"
$ws.Range("E29").Value = "  -1.80%  "
$ws.Range("E30").Value = "  -0.41%  "
$ws.Range("E31").Value = "  -1.11%  "
$ws.Range("D32").Value = "3.35"
$ws.Range("E32").Value = "  -0.09%  "
$ws.Range("D33").Value = "2.98"
$ws.Range("E33").Value = "  -1.81%  "
$ws.Range("D34").Value = "1.309.36"
$ws.Range("E34").Value = "  +2.91%  "
$ws.Range("D35").Value = "1.55"
$ws.Range("E35").Value = "  -1.75%  "
$ws.Range("E36").Value = "  -0.18%  "
$ws.Range("D37").Value = "0.0174"
$ws.Range("E37").Value = "  -2.16%  "
$ws.Range("D38").Value = "0.537"
$ws.Range("E38").Value = "  -1.25%  "
$ws.Range("D39").Value = "0.839"
$ws.Range("E39").Value = "  -0.93%  "
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("D41").Value = "2.26"
$ws.Range("E41").Value = "  +1.81%  "
$ws.Range("D42").Value = "0.803"
$ws.Range("E42").Value = "  -1.00%  "
$ws.Range("D43").Value = "5.25"
$ws.Range("E43").Value = "  -2.39%  "
$ws.Range("D44").Value = "1.766.20"
$ws.Range("E44").Value = "  -1.30%  "
$ws.Range("D45").Value = "62.39"
$ws.Range("E45").Value = "  -0.34%  "
$ws.Range("D46").Value = "90.63"
$ws.Range("E46").Value = "  -1.81%  "
$ws.Range("D47").Value = "1.59"
$ws.Range("E47").Value = "  -0.31%  "
$ws.Range("E48").Value = "  +17.41%  "
$ws.Range("D49").Value = "0.797"
$ws.Range("E49").Value = "  +17.71%  "
$ws.Range("E50").Value = "  -0.16%  "
$ws.Range("D51").Value = "7.54"
$ws.Range("E51").Value = "  -2.48%  "
